# Updated symbol list on Wed Feb  8 04:37:00 UTC 2023 with GitHub Actions
#
# Re-applies one polling cycle of the coinranking.com scraper: the coin at
# each rank (rows 7-18) rotates down one slot (today's #N becomes
# tomorrow's #N+1, wrapping GateToken back to the top of that block) and
# every Price/Volume(1h) cell across the sheet is refreshed with the
# latest reading. Price/Volume cells are stored as plain text in this
# workbook (e.g. "332.14", "1.52%"), so each write forces the Text number
# format first and clears it back off afterwards -- otherwise COM would
# silently coerce the numeric-looking strings into real numbers/percentages
# and drop the original text formatting (trailing zeros, "%" suffix, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

# Coin / Link (plain text -- no numeric coercion risk)
Set-TextValue "B7"  "MXToken"
Set-TextValue "C7"  "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "B8"  "BTSEToken"
Set-TextValue "C8"  "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "B9"  "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9"  "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "B11" "MCDex"
Set-TextValue "C11" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "B18" "GateToken"
Set-TextValue "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

# Price / Volume(1h) -- numeric-looking text, needs the Text-format guard
Set-TextValue "D2"  "332.14"
Set-TextValue "E2"  "1.52%"
Set-TextValue "D3"  "45.86"
Set-TextValue "E3"  "4.44%"
Set-TextValue "D4"  "5.665"
Set-TextValue "E4"  "2.73%"
Set-TextValue "D5"  "0.08374"
Set-TextValue "E5"  "4.39%"
Set-TextValue "D6"  "2.037"
Set-TextValue "E6"  "2.52%"
Set-TextValue "D7"  "0.9847"
Set-TextValue "E7"  "3.31%"
Set-TextValue "D8"  "2.586"
Set-TextValue "E8"  "0.45%"
Set-TextValue "D9"  "0.1153"
Set-TextValue "E9"  "2.51%"
Set-TextValue "D10" "0.1938"
Set-TextValue "E10" "3.72%"
Set-TextValue "D11" "10.39"
Set-TextValue "E11" "-2.84%"
Set-TextValue "D12" "0.1006"
Set-TextValue "E12" "2.43%"
Set-TextValue "D13" "0.04663"
Set-TextValue "E13" "1.70%"
Set-TextValue "D14" "0.1058"
Set-TextValue "E14" "-0.87%"
Set-TextValue "D15" "0.001290"
Set-TextValue "E15" "1.44%"
Set-TextValue "D16" "0.006056"
Set-TextValue "E16" "3.92%"
Set-TextValue "D17" "3.368"
Set-TextValue "E17" "0.43%"
Set-TextValue "D18" "4.488"
Set-TextValue "E18" "4.49%"

Set-TextValue "D20" "0.1399"
Set-TextValue "E20" "-0.57%"
Set-TextValue "D21" "0.2606"
Set-TextValue "E21" "2.38%"
Set-TextValue "D22" "0.04218"
Set-TextValue "E22" "3.28%"
Set-TextValue "D23" "0.001309"
Set-TextValue "E23" "5.15%"
Set-TextValue "D24" "0.004660"
Set-TextValue "E24" "7.53%"
Set-TextValue "D25" "0.0001281"
Set-TextValue "E25" "10.62%"
Set-TextValue "D26" "0.0003740"
Set-TextValue "E26" "-0.09%"

Set-TextValue "D38" "0.02802"
Set-TextValue "E38" "9.71%"
Set-TextValue "D39" "0.05808"
Set-TextValue "E39" "2.44%"
Set-TextValue "D40" "0.007739"
Set-TextValue "E40" "2.50%"
Set-TextValue "D41" "0.1438"
Set-TextValue "E41" "2.68%"
Set-TextValue "D42" "0.007190"
Set-TextValue "E42" "-5.31%"
Set-TextValue "D43" "0.001974"
Set-TextValue "E43" "-1.91%"
Set-TextValue "D44" "0.008067"
Set-TextValue "E44" "-5.24%"
Set-TextValue "D45" "0.3499"
Set-TextValue "D46" "0.00007299"
Set-TextValue "E46" "2.63%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.02%"
Set-TextValue "D48" "0.0005801"
Set-TextValue "E48" "-0.19%"
Set-TextValue "D49" "0.003498"
Set-TextValue "E49" "13.35%"
Set-TextValue "D50" "0.003498"
Set-TextValue "E50" "-0.87%"
Set-TextValue "D51" "0.00002100"
Set-TextValue "E51" "0.02%"
